# Auto-generated Excel COM-interop edit script
# Applies numeric updates to H:N columns across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3903.027
$ws.Range("I43").Value = 2588.3333
$ws.Range("J43").Value = 4157.484
$ws.Range("K43").Value = 2588.3333
$ws.Range("L43").Value = 4157.484
$ws.Range("M43").Value = -2519.3333
$ws.Range("N43").Value = -4295.484

$ws.Range("H132").Value = 4165.3584
$ws.Range("I132").Value = 3672.7346
$ws.Range("J132").Value = 10200
$ws.Range("K132").Value = 11018.2038
$ws.Range("L132").Value = 30600
$ws.Range("M132").Value = -8488.203799999999
$ws.Range("N132").Value = -35660

$ws.Range("H141").Value = 1138.6842
$ws.Range("I141").Value = 850.30304
$ws.Range("K141").Value = 2550.90912
$ws.Range("M141").Value = 2629.09088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1375
$ws.Range("I63").Value = 1500
$ws.Range("J63").Value = 1000
$ws.Range("K63").Value = 1500
$ws.Range("L63").Value = 1000
$ws.Range("M63").Value = -814
$ws.Range("N63").Value = -2372

$ws.Range("H66").Value = 1375
$ws.Range("I66").Value = 1500
$ws.Range("J66").Value = 1000
$ws.Range("K66").Value = 7500
$ws.Range("L66").Value = 5000
$ws.Range("M66").Value = -4068
$ws.Range("N66").Value = -11864

$ws.Range("H74").Value = 21742904
$ws.Range("I74").Value = 35716424
$ws.Range("J74").Value = 6314.222
$ws.Range("K74").Value = 35716424
$ws.Range("L74").Value = 6314.222
$ws.Range("M74").Value = -35715550
$ws.Range("N74").Value = -8062.222

$ws.Range("H77").Value = 21742904
$ws.Range("I77").Value = 35716424
$ws.Range("J77").Value = 6314.222
$ws.Range("K77").Value = 178582120
$ws.Range("L77").Value = 31571.11
$ws.Range("M77").Value = -178577752
$ws.Range("N77").Value = -40307.11

$ws.Range("H88").Value = 2930.9333
$ws.Range("I88").Value = 3108.3333
$ws.Range("J88").Value = 2812.6667
$ws.Range("K88").Value = 3108.3333
$ws.Range("L88").Value = 2812.6667
$ws.Range("M88").Value = -2702.3333
$ws.Range("N88").Value = -3624.6667

$ws.Range("H91").Value = 2930.9333
$ws.Range("I91").Value = 3108.3333
$ws.Range("J91").Value = 2812.6667
$ws.Range("K91").Value = 3108.3333
$ws.Range("L91").Value = 2812.6667
$ws.Range("M91").Value = -1704.3333
$ws.Range("N91").Value = -5620.6667

$ws.Range("H132").Value = 5557539.5
$ws.Range("I132").Value = 7354533
$ws.Range("J132").Value = 3196
$ws.Range("K132").Value = 22063599
$ws.Range("L132").Value = 9588
$ws.Range("M132").Value = -22061069
$ws.Range("N132").Value = -14648

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1968.2678
$ws.Range("I134").Value = 1246.8605
$ws.Range("J134").Value = 4354.4614
$ws.Range("K134").Value = 3740.5815
$ws.Range("L134").Value = 13063.3842
$ws.Range("M134").Value = -1205.5815
$ws.Range("N134").Value = -18133.3842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 13515822
$ws.Range("I132").Value = 23811748
$ws.Range("J132").Value = 2420.625
$ws.Range("K132").Value = 71435244
$ws.Range("L132").Value = 7261.875
$ws.Range("M132").Value = -71432714
$ws.Range("N132").Value = -12321.875

$ws.Range("H134").Value = 1489.125
$ws.Range("I134").Value = 1473.3214
$ws.Range("J134").Value = 1526
$ws.Range("K134").Value = 4419.9642
$ws.Range("L134").Value = 4578
$ws.Range("M134").Value = -1884.9642
$ws.Range("N134").Value = -9648

$ws.Range("H140").Value = 46465.57
$ws.Range("J140").Value = 46465.57
$ws.Range("L140").Value = 46465.57
$ws.Range("N140").Value = -56825.57

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 225.5
$ws.Range("I17").Value = 100
$ws.Range("J17").Value = 351
$ws.Range("K17").Value = 300
$ws.Range("L17").Value = 1053
$ws.Range("M17").Value = -131
$ws.Range("N17").Value = -1391

$ws.Range("H23").Value = 117.77273
$ws.Range("J23").Value = 136.66667
$ws.Range("L23").Value = 410.00001
$ws.Range("N23").Value = -880.00001

$ws.Range("H55").Value = 523.6842
$ws.Range("I55").Value = 200
$ws.Range("J55").Value = 610
$ws.Range("K55").Value = 600
$ws.Range("L55").Value = 1830
$ws.Range("M55").Value = -423
$ws.Range("N55").Value = -2184

$ws.Range("H68").Value = 740
$ws.Range("I68").Value = 570.25
$ws.Range("J68").Value = 824.875
$ws.Range("K68").Value = 1710.75
$ws.Range("L68").Value = 2474.625
$ws.Range("M68").Value = -899.75
$ws.Range("N68").Value = -4096.625

$ws.Range("H71").Value = 740
$ws.Range("I71").Value = 570.25
$ws.Range("J71").Value = 824.875
$ws.Range("K71").Value = 5132.25
$ws.Range("L71").Value = 7423.875
$ws.Range("M71").Value = -1076.25
$ws.Range("N71").Value = -15535.875

$ws.Range("H94").Value = 6053
$ws.Range("I94").Value = 2712
$ws.Range("J94").Value = 7166.6665
$ws.Range("K94").Value = 8136
$ws.Range("L94").Value = 21499.9995
$ws.Range("M94").Value = -7460
$ws.Range("N94").Value = -22851.9995

$ws.Range("H98").Value = 914.3889
$ws.Range("I98").Value = 465
$ws.Range("J98").Value = 1139.0834
$ws.Range("K98").Value = 1395
$ws.Range("L98").Value = 3417.2502
$ws.Range("M98").Value = 103
$ws.Range("N98").Value = -6413.2502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 230000
$ws.Range("I53").Value = 230000
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 230000
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -229369
$ws.Range("N53").ClearContents()

$ws.Range("H132").Value = 2435.8494
$ws.Range("I132").Value = 1840.0392
$ws.Range("J132").Value = 3817.0454
$ws.Range("K132").Value = 5520.1176
$ws.Range("L132").Value = 11451.1362
$ws.Range("M132").Value = -2990.1176
$ws.Range("N132").Value = -16511.1362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5078.394
$ws.Range("I7").Value = 4982.638
$ws.Range("J7").Value = 5315.263
$ws.Range("K7").Value = 4982.638
$ws.Range("L7").Value = 5315.263
$ws.Range("M7").Value = -4870.638
$ws.Range("N7").Value = -5539.263

$ws.Range("H122").Value = 4639.4287
$ws.Range("I122").Value = 5172
$ws.Range("J122").Value = 3816.3635
$ws.Range("K122").Value = 15516
$ws.Range("L122").Value = 11449.0905
$ws.Range("M122").Value = -13066
$ws.Range("N122").Value = -16349.0905

$ws.Range("H126").Value = 5078.394
$ws.Range("I126").Value = 4982.638
$ws.Range("J126").Value = 5315.263
$ws.Range("K126").Value = 14947.914
$ws.Range("L126").Value = 15945.789
$ws.Range("M126").Value = -12477.914
$ws.Range("N126").Value = -20885.789

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1685.0364
$ws.Range("I132").Value = 1423.4318
$ws.Range("J132").Value = 2731.4546
$ws.Range("K132").Value = 4270.2954
$ws.Range("L132").Value = 8194.363799999999
$ws.Range("M132").Value = -1740.2954
$ws.Range("N132").Value = -13254.3638
